# Apply the "address.city" column addition to the req-payload sheet and
# adjust the active sheet/selection bookkeeping to match.

$wb = $excel.ActiveWorkbook

$wsReq = $wb.Worksheets.Item("req-payload")

# --- req-payload sheet: insert a new row (row 3) so id=1 gets a second
#     row holding the nested "address.city" value, and add a new column D
#     ("address.city") with Delhi/Bombay values. ---

# Shift the old row 3 (id=2, dorpheus, follower) down to row 4.
$wsReq.Rows.Item(3).Insert()

# New header + values for column D (added in top-to-bottom order so the
# shared-strings table picks up "address.city", "Delhi", "Bombay" in that
# order).
$wsReq.Range("D1").Value = "address.city"
$wsReq.Range("D2").Value = "Delhi"

# New row 3: id repeats (1) and only the new D column is populated.
$wsReq.Range("A3").Value = 1
$wsReq.Range("D3").Value = "Bombay"

# Column widths: Excel auto-fits B:D to the new content (recorded in the
# saved file as bestFit custom widths of 10 / 8.57 / 11.57 characters).
$wsReq.Columns.Item(2).ColumnWidth = 9.166666666666666
$wsReq.Columns.Item(3).ColumnWidth = 7.666666666666667
$wsReq.Columns.Item(4).ColumnWidth = 10.666666666666666

# Selection / active-cell bookkeeping on req-payload.
$wsReq.Range("X4").Select()

# --- sheet activation: req-payload becomes the active/selected tab,
#     json-path loses tabSelected. ---
$wsReq.Activate()

$wb.Save()
